# Update simulation results (gw_simul, gcpi_simul, diffcpicf_simul, cf10_simul, cf1_simul)
# for rows 6-19 with refreshed dynamic simulation output, in support of new graphing functionality.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updatedValues = @{
    6 = @{ "C" = 2.920844475933321; "D" = 3.50332108551396; "E" = 0.3651035214814868; "F" = 1.629437769201222; "G" = 1.74575349881511 }
    7 = @{ "C" = 2.659297480106434; "D" = -0.8097245042667431; "E" = 0.2023071000046448; "F" = 1.645171086450825; "G" = 1.218188589753771 }
    8 = @{ "C" = 2.356016008563889; "D" = 0.6587158627226681; "E" = -0.2425954644088348; "F" = 1.228150973030189; "G" = 0.8777762550840216 }
    9 = @{ "C" = -0.6720430974772955; "D" = 3.808187847745355; "E" = 0.1451970944131848; "F" = 1.408362456742319; "G" = 1.38042602442851 }
    10 = @{ "C" = 1.387983868786141; "D" = 2.281759681484334; "E" = -0.2610187768937071; "F" = 1.387084963369895; "G" = 1.641575190993414 }
    11 = @{ "C" = 4.273672589567485; "D" = 5.484468090056484; "E" = 1.840094280748439; "F" = 1.489193676687087; "G" = 1.796488075003597 }
    12 = @{ "C" = 2.235732564447486; "D" = 6.361241609291474; "E" = 3.606138052060389; "F" = 1.652571091429816; "G" = 2.293974906802815 }
    13 = @{ "C" = 3.225634172325693; "D" = 3.397113653820869; "E" = 3.00071973423478; "F" = 1.62125614033817; "G" = 2.383909561241147 }
    14 = @{ "C" = 4.233872935525688; "D" = 7.611732589516219; "E" = 4.072063794677848; "F" = 1.7534828537285; "G" = 2.372833234470688 }
    15 = @{ "C" = 4.501937528334725; "D" = 7.823379951184652; "E" = 4.501878875949707; "F" = 1.948100491367711; "G" = 2.97965742450988 }
    16 = @{ "C" = 4.666060466001984; "D" = 8.235016640234173; "E" = 4.472835801886164; "F" = 2.071549268806558; "G" = 3.668721205996093 }
    17 = @{ "C" = 4.425877292785612; "D" = 0.9032026041690222; "E" = 3.759423385034869; "F" = 2.260981562323539; "G" = 3.358674458632915 }
    18 = @{ "C" = 4.062454202254666; "D" = 3.675076561882506; "E" = 2.786335704896901; "F" = 2.241176634204947; "G" = 2.28953085825381 }
    19 = @{ "C" = 4.110657281150679; "D" = 3.297099956269078; "E" = 1.047941516128815; "F" = 2.30947278585226; "G" = 2.592748903861609 }
}

foreach ($rowNum in $updatedValues.Keys) {
    $rowValues = $updatedValues[$rowNum]
    foreach ($col in $rowValues.Keys) {
        $cellRef = "$col$rowNum"
        $ws.Range($cellRef).Value = $rowValues[$col]
    }
}
